# Auto-generated-by-analysis PowerShell COM-interop script
# Applies updated market-board snapshot values (currentAveragePrice*, LevePrice*, LeveProfit*)
# to the per-job Leve profitability sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5475.0205
$ws.Range("J17").Value = 5475.0205
$ws.Range("L17").Value = 16425.0615
$ws.Range("N17").Value = -16761.0615
$ws.Range("H70").Value = 1761.7333
$ws.Range("I70").Value = 1833.5
$ws.Range("J70").Value = 1713.8889
$ws.Range("K70").Value = 5500.5
$ws.Range("L70").Value = 5141.6667
$ws.Range("M70").Value = -5230.5
$ws.Range("N70").Value = -5681.6667
$ws.Range("H73").Value = 1761.7333
$ws.Range("I73").Value = 1833.5
$ws.Range("J73").Value = 1713.8889
$ws.Range("K73").Value = 5500.5
$ws.Range("L73").Value = 5141.6667
$ws.Range("M73").Value = -4564.5
$ws.Range("N73").Value = -7013.6667
$ws.Range("H86").Value = 8279.5
$ws.Range("I86").Value = 7002.5835
$ws.Range("K86").Value = 7002.5835
$ws.Range("M86").Value = -5879.5835
$ws.Range("H89").Value = 8279.5
$ws.Range("I89").Value = 7002.5835
$ws.Range("K89").Value = 35012.9175
$ws.Range("M89").Value = -29396.9175
$ws.Range("H112").Value = 2116.5833
$ws.Range("J112").Value = 2172.6365
$ws.Range("L112").Value = 6517.9095
$ws.Range("N112").Value = -8733.9095
$ws.Range("H132").Value = 2073.6406
$ws.Range("I132").Value = 2223.111
$ws.Range("J132").Value = 1266.5
$ws.Range("K132").Value = 6669.333
$ws.Range("L132").Value = 3799.5
$ws.Range("M132").Value = -4139.333
$ws.Range("N132").Value = -8859.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("H11").Value = 1672.6666
$ws.Range("J11").Value = 2507.5
$ws.Range("L11").Value = 2507.5
$ws.Range("N11").Value = -2795.5
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = ""
$ws.Range("H32").Value = 11135.811
$ws.Range("I32").Value = 7559.738
$ws.Range("K32").Value = 7559.738
$ws.Range("M32").Value = -7272.738
$ws.Range("H74").Value = 27623.674
$ws.Range("I74").Value = 6865.3613
$ws.Range("K74").Value = 6865.3613
$ws.Range("M74").Value = -5991.3613
$ws.Range("H77").Value = 27623.674
$ws.Range("I77").Value = 6865.3613
$ws.Range("K77").Value = 34326.8065
$ws.Range("M77").Value = -29958.8065
$ws.Range("H97").Value = 3357021.8
$ws.Range("I97").Value = 4195270.5
$ws.Range("K97").Value = 4195270.5
$ws.Range("M97").Value = -4194774.5
$ws.Range("H102").Value = 5560193.5
$ws.Range("J102").Value = 4166.6665
$ws.Range("L102").Value = 4166.6665
$ws.Range("N102").Value = -7410.6665
$ws.Range("H122").Value = 3502664.5
$ws.Range("I122").Value = 6581208
$ws.Range("K122").Value = 19743624
$ws.Range("M122").Value = -19741174
$ws.Range("H132").Value = 26392.62
$ws.Range("I132").Value = 5993.364
$ws.Range("K132").Value = 17980.092
$ws.Range("M132").Value = -15450.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2442197.8
$ws.Range("I94").Value = 2632739.8
$ws.Range("J94").Value = 28666.666
$ws.Range("K94").Value = 2632739.8
$ws.Range("L94").Value = 28666.666
$ws.Range("M94").Value = -2632288.8
$ws.Range("N94").Value = -29568.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 16992.5
$ws.Range("J9").Value = 16992.5
$ws.Range("L9").Value = 16992.5
$ws.Range("N9").Value = -17328.5
$ws.Range("H36").Value = 9995.666999999999
$ws.Range("I36").Value = 9995.666999999999
$ws.Range("K36").Value = 9995.666999999999
$ws.Range("M36").Value = -9607.666999999999
$ws.Range("H40").Value = 9995.666999999999
$ws.Range("I40").Value = 9995.666999999999
$ws.Range("K40").Value = 9995.666999999999
$ws.Range("M40").Value = -9835.666999999999
$ws.Range("H58").Value = 5933.24
$ws.Range("I58").Value = 7776.8125
$ws.Range("K58").Value = 7776.8125
$ws.Range("M58").Value = -7573.8125
$ws.Range("H86").Value = 6925.8076
$ws.Range("I86").Value = 5158.3
$ws.Range("K86").Value = 5158.3
$ws.Range("M86").Value = -4035.3
$ws.Range("H89").Value = 6925.8076
$ws.Range("I89").Value = 5158.3
$ws.Range("K89").Value = 25791.5
$ws.Range("M89").Value = -20175.5
$ws.Range("H122").Value = 1809.25
$ws.Range("I122").Value = 2054.5
$ws.Range("K122").Value = 6163.5
$ws.Range("M122").Value = -3713.5
$ws.Range("H136").Value = 5933.24
$ws.Range("I136").Value = 7776.8125
$ws.Range("K136").Value = 23330.4375
$ws.Range("M136").Value = -20780.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 38536.89
$ws.Range("I5").Value = 995.1111
$ws.Range("J5").Value = 113620.445
$ws.Range("K5").Value = 2985.3333
$ws.Range("L5").Value = 340861.335
$ws.Range("M5").Value = -2873.3333
$ws.Range("N5").Value = -341085.335
$ws.Range("H56").Value = 10421766
$ws.Range("I56").Value = 10421766
$ws.Range("K56").Value = 10421766
$ws.Range("M56").Value = -10421236
$ws.Range("H135").Value = 38536.89
$ws.Range("I135").Value = 995.1111
$ws.Range("J135").Value = 113620.445
$ws.Range("K135").Value = 8955.999899999999
$ws.Range("L135").Value = 1022584.005
$ws.Range("M135").Value = -6420.999899999999
$ws.Range("N135").Value = -1027654.005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""
$ws.Range("H15").Value = 29797
$ws.Range("J15").Value = 29797
$ws.Range("L15").Value = 29797
$ws.Range("N15").Value = -30373
$ws.Range("H81").Value = 29797
$ws.Range("J81").Value = 29797
$ws.Range("L81").Value = 29797
$ws.Range("N81").Value = -31793
$ws.Range("H84").Value = 29797
$ws.Range("J84").Value = 29797
$ws.Range("L84").Value = 89391
$ws.Range("N84").Value = -99375
$ws.Range("H92").Value = 18125.5
$ws.Range("J92").Value = 18125.5
$ws.Range("L92").Value = 18125.5
$ws.Range("N92").Value = -21869.5
$ws.Range("H97").Value = 1323614.8
$ws.Range("I97").Value = 1984944.8
$ws.Range("J97").Value = 954.8333
$ws.Range("K97").Value = 1984944.8
$ws.Range("L97").Value = 954.8333
$ws.Range("M97").Value = -1984448.8
$ws.Range("N97").Value = -1946.8333
$ws.Range("H119").Value = 95695
$ws.Range("J119").Value = 95695
$ws.Range("L119").Value = 95695
$ws.Range("N119").Value = -105371
$ws.Range("H122").Value = 406766.38
$ws.Range("I122").Value = 687113.6
$ws.Range("J122").Value = 1820.3334
$ws.Range("K122").Value = 2061340.8
$ws.Range("L122").Value = 5461.0002
$ws.Range("M122").Value = -2058890.8
$ws.Range("N122").Value = -10361.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5123.0835
$ws.Range("I96").Value = 4996.25
$ws.Range("J96").Value = 5186.5
$ws.Range("K96").Value = 4996.25
$ws.Range("L96").Value = 5186.5
$ws.Range("M96").Value = -3623.25
$ws.Range("N96").Value = -7932.5
$ws.Range("H100").Value = 2050.2354
$ws.Range("I100").Value = 2466.6155
$ws.Range("J100").Value = 697
$ws.Range("K100").Value = 4933.231
$ws.Range("L100").Value = 1394
$ws.Range("M100").Value = -4392.231
$ws.Range("N100").Value = -2476
$ws.Range("H122").Value = 3464.4827
$ws.Range("I122").Value = 1932.5
$ws.Range("K122").Value = 5797.5
$ws.Range("M122").Value = -3347.5
$ws.Range("H126").Value = 2095.739
$ws.Range("I126").Value = 2008.1666
$ws.Range("K126").Value = 6024.4998
$ws.Range("M126").Value = -3554.4998

